$d = $word.ActiveDocument
$br = [char]11

$old1 = "The core idea behind the self-attention mechanism, as proposed in the paper `"Attention Is All You Need,`" is to compute a representation of a sequence by relating different positions within that sequence. Self-attention, also known as intra-attention, allows the model to weigh the importance of each position in the sequence relative to others, enabling it to capture dependencies between distant positions efficiently."
$new1 = "The core idea behind the self-attention mechanism, as proposed in the paper `"Attention Is All You Need,`" is to compute a representation of a sequence by relating different positions within that sequence. This mechanism, also known as intra-attention, allows the model to focus on different parts of the input sequence when constructing its output representation."
$rng1 = $d.Content
$found1 = $rng1.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false)
if ($found1) {
    $rng1.Text = $new1
}
Write-Output "Replace 1: $found1"

$old2 = "In traditional sequence models, capturing long-range dependencies can be computationally expensive and challenging. However, self-attention reduces this complexity to a constant number of operations, allowing the model to focus on relevant parts of the sequence regardless of their distance from each other. This is achieved by calculating attention scores for each position, which are then used to create a weighted sum of the input features, effectively highlighting important parts of the sequence."
$new2 = "Self-attention operates by assigning attention scores to each position in the sequence relative to every other position. These scores determine how much influence each position should have in the final representation of a particular position. The attention scores are typically computed using a compatibility function, such as the dot product, followed by a softmax operation to ensure they sum to one."
$rng2 = $d.Content
$found2 = $rng2.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false)
if ($found2) {
    $rng2.Text = $new2
}
Write-Output "Replace 2: $found2"

$old3 = "The paper addresses a potential downside of self-attention, which is the reduced effective resolution due to averaging attention-weighted positions. This issue is mitigated through the use of Multi-Head Attention, which allows the model to attend to information from different representation subspaces at different positions, thereby enhancing its ability to capture complex patterns and relationships within the data."
$new3 = "One of the key advantages of self-attention is its ability to capture dependencies between distant positions in the sequence efficiently. In traditional sequence models, such as recurrent neural networks (RNNs), capturing long-range dependencies can be computationally expensive and challenging. In contrast, self-attention reduces this complexity to a constant number of operations, regardless of the distance between positions.${br}${br}However, this efficiency comes at the cost of reduced effective resolution due to the averaging of attention-weighted positions. To address this, the paper introduces Multi-Head Attention, which allows the model to attend to information from multiple representation subspaces simultaneously, thereby enhancing the model's ability to capture complex dependencies."
$rng3 = $d.Content
$found3 = $rng3.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false)
if ($found3) {
    $rng3.Text = $new3
}
Write-Output "Replace 3: $found3"

$old4 = "Self-attention has been successfully applied to various tasks, including reading comprehension, abstractive summarization, textual entailment, and learning task-independent sentence representations, demonstrating its versatility and effectiveness in processing sequential data."
$new4 = "Self-attention has been successfully applied to various tasks, including reading comprehension, abstractive summarization, textual entailment, and learning task-independent sentence representations, demonstrating its versatility and effectiveness in natural language processing."
$rng4 = $d.Content
$found4 = $rng4.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false)
if ($found4) {
    $rng4.Text = $new4
}
Write-Output "Replace 4: $found4"

Write-Output $d.Content.Text